$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 4 (even_MAG-GUT70200.fa) - entire row removed
$ws.Rows.Item(4).Delete()

# Delete column C ("max") - shifts former D (prediction) and E (rejection-f) left
$ws.Columns.Item(3).Delete()

# Update the numeric values in column B for remaining data rows
$ws.Cells.Item(2, 2).Value = 4.52389648757414
$ws.Cells.Item(3, 2).Value = 2.927948143320764
